$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Duplicate the "while" test-case block (cols C:F) into new cols G:J,
#    carrying over cell formatting first, then the (mostly identical) values.
# ---------------------------------------------------------------------------
$ws.Range("C1:F3").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

# Header row (row 1) - G:J mirror C:F exactly
$ws.Range("G1").Value = "while"
$ws.Range("H1").Value = "setProperty"
$ws.Range("I1").Value = "echo"
$ws.Range("J1").Value = "end"

# Row 2 - G:J mirror C:F (condition cell G2 left blank, like F2)
$ws.Range("H2").Value = "json"
$ws.Range("I2").Value = '${userName}_${times}'

# Row 3 - new "target" json replaces old condition string, D3/H3 stay the same "times" increment
$ws.Range("G3").Value = '{"target":"userName == ''hugang'' && parseInt(times) <= 10"}'
$ws.Range("H3").Value = '{"times":${times+1}}'

# ---------------------------------------------------------------------------
# 2. Original columns: E2 now carries the "${userName}_${times}" value while
#    E3's old condition string moves out to the new G3 cell (so E3 clears).
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = '${userName}_${times}'
$ws.Range("E3").ClearContents()

# ---------------------------------------------------------------------------
# 3. Column widths (best effort through the COM ColumnWidth property - the
#    host quantizes to whole pixels, so these are the closest achievable
#    values to the target raw widths of 26.5 / 57.375 / 20.5 / 22.375 / 4.75).
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 25.714285714285715    # col E: -> 26.5
$ws.Columns.Item(7).ColumnWidth = 56.714285714285715    # col G: -> 57.375
$ws.Columns.Item(8).ColumnWidth = 19.714285714285715    # col H: -> 20.5
$ws.Columns.Item(9).ColumnWidth = 21.714285714285715    # col I: -> 22.375
$ws.Columns.Item(10).ColumnWidth = 4.0                  # col J: -> 4.75

# ---------------------------------------------------------------------------
# 4. Selection / view state.
# ---------------------------------------------------------------------------
$ws.Range("G17").Select() | Out-Null
